# Insert a new data row at row 526 ("Vega Modelo de Temuco" / Mango price
# list). This shifts every existing row from 526..638 down to 527..639,
# matching the recorded edit (dimension grows from A1:T638 to A1:T639).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(526).Insert()

$ws.Range("A526").Value = 10
$ws.Range("B526").Value = "Vega Modelo de Temuco"
$ws.Range("C526").Value = "La Araucanía"
$ws.Range("D526").Value = 45173
$ws.Range("E526").Value = 9
$ws.Range("F526").Value = "Fruta"
$ws.Range("G526").Value = 100108
$ws.Range("H526").Value = "Tropicales y subtropicales"
$ws.Range("I526").Value = 100108002
$ws.Range("J526").Value = "Mango"
$ws.Range("K526").Value = "Sin especificar"
$ws.Range("L526").Value = "Primera"
$ws.Range("M526").Value = 800
$ws.Range("N526").Value = 10000
$ws.Range("O526").Value = 12000
$ws.Range("P526").Value = 11250
$ws.Range("Q526").Value = '$/bandeja 4 kilos'
$ws.Range("R526").Value = "Brasil"
$ws.Range("S526").Value = 2812
$ws.Range("T526").Value = 4
